# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp in A1
# - Swap the Tailandia / Tayikistan rows (A79 <-> A80 labels) and refresh
#   their case figures
# - Refresh case-count figures (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for the other
#   updated countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 18:05"

# --- Tailandia / Tayikistan swap (row 79 / row 80) ---------------------
# Before: A79 = Tailandia, A80 = Tayikistan
# After:  A79 = Tayikistan, A80 = Tailandia (labels swap, data refreshed)
$ws.Range("A79").Value = "Tayikistan"
$ws.Range("B79:H79").Value = @(@(3100, 171, 1395, 1659, 0, 0, 46))

$ws.Range("A80").Value = "Tailandia"
$ws.Range("B80:H80").Value = @(@(3042, 2, 2928, 57, 0, 1, 57))

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4:H4").Value = @(@(1691612, 5176, 451857, 1140347, 0, 108, 99408))

# --- Alemania (row 11) --------------------------------------------------
$ws.Range("B11:H11").Value = @(@(180566, 238, 161200, 10977, 0, 18, 8389))

# --- India (row 13) ------------------------------------------------------
$ws.Range("B13:E13").Value = @(@(144118, 5582, 59913, 80088))

# --- Canada (row 16) ------------------------------------------------------
$ws.Range("B16:H16").Value = @(@(85103, 404, 44206, 34444, 0, 29, 6453))

# --- Polonia (row 37) ------------------------------------------------------
$ws.Range("B37:H37").Value = @(@(21631, 305, 9276, 11348, 0, 11, 1007))

# --- Chequia (row 54) ------------------------------------------------------
$ws.Range("B54:H54").Value = @(@(8972, 17, 6180, 2475, 0, 2, 317))

# --- Luxemburgo (row 71) ------------------------------------------------------
$ws.Range("B71:E71").Value = @(@(3993, 1, 3767, 116))

# --- Republica de Chipre (row 115) ------------------------------------------------------
$ws.Range("B115:E115").Value = @(@(937, 2, 594, 326))

# --- Jordania (row 125) ------------------------------------------------------
$ws.Range("B125:E125").Value = @(@(711, 3, 479, 223))
